# The author fixed a typo in the bandwidth table: cell E3 ("Remote Disk" row
# for the "BW" series) was off by a decimal place (18.034764 -> 180.34764).
# Update the cell value, then leave the selection where Excel would land
# after typing the value into E3 and pressing Enter (i.e. E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 180.34764000000001

[void]$ws.Range("E4").Select()
